$d = $word.ActiveDocument
$t1 = $d.Tables.Item(1)

# --- 1. Split "Player name" (row 3, col 2) into two runs: "First " and "name" ---
$cell = $t1.Cell(3, 2)
$cell.Range.Text = "First name"

# Force a run boundary after "First " (6 chars incl. trailing space) by
# toggling a character property on and back off; this splits the single
# run into two runs without altering the visible formatting.
$cellAfterText = $t1.Cell(3, 2)
$splitRange = $d.Range($cellAfterText.Range.Start, $cellAfterText.Range.Start + 6)
$splitRange.Bold = 1
$cellAfterBold = $t1.Cell(3, 2)
$splitRangeBack = $d.Range($cellAfterBold.Range.Start, $cellAfterBold.Range.Start + 6)
$splitRangeBack.Bold = 0

# --- 2. Add a new row to the table: "String" | "Last name" (with _GoBack bookmark) ---
$newRow = $t1.Rows.Add()

$newCell1 = $t1.Cell(4, 1)
$newCell1.Range.Text = "String"

$newCell2 = $t1.Cell(4, 2)
$newCell2.Range.Text = "Last nameZ"

# Insert the _GoBack bookmark right after "Last name" (before the temporary
# placeholder character). A bookmark can't be collapsed exactly at the end
# of a paragraph's content in one step, so we anchor it just before a
# placeholder character and then delete that placeholder.
$cellForBm = $t1.Cell(4, 2)
$bmPos = $cellForBm.Range.Start + 9
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cellWithPlaceholder = $t1.Cell(4, 2)
$placeholderRange = $d.Range($cellWithPlaceholder.Range.End - 2, $cellWithPlaceholder.Range.End - 1)
$placeholderRange.Delete()

# (Adding the "_GoBack" bookmark above automatically removed the previous
# "_GoBack" bookmark that sat at the end of the "Create a junction table..."
# paragraph, matching the target diff.)
